# Auto-generated edit script applying scheduled-runner value updates
# to the Tonberry Profits workbook (per-sheet leve profit recalculation).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 11
$ws.Range("H11").Value = 15384748
$ws.Range("I11").Value = 15384748
$ws.Range("K11").Value = 15384748
$ws.Range("M11").Value = -15384608
# Row 15
$ws.Range("H15").Value = 803.10254
$ws.Range("I15").Value = 803.10254
$ws.Range("K15").Value = 2409.30762
$ws.Range("M15").Value = -2240.30762
# Row 80
$ws.Range("H80").Value = 1528.2667
$ws.Range("I80").Value = 1551.8
$ws.Range("K80").Value = 4655.4
$ws.Range("M80").Value = -3657.4
# Row 83
$ws.Range("H83").Value = 1528.2667
$ws.Range("I83").Value = 1551.8
$ws.Range("K83").Value = 13966.2
$ws.Range("M83").Value = -8974.199999999999
# Row 132
$ws.Range("H132").Value = 1160.4062
$ws.Range("I132").Value = 1146.6786
$ws.Range("J132").Value = 1256.5
$ws.Range("K132").Value = 3440.0358
$ws.Range("L132").Value = 3769.5
$ws.Range("M132").Value = -910.0357999999997
$ws.Range("N132").Value = -8829.5
# Row 138
$ws.Range("H138").Value = 2154.898
$ws.Range("I138").Value = 2308.1538
$ws.Range("J138").Value = 2053.5933
$ws.Range("K138").Value = 6924.4614
$ws.Range("L138").Value = 6160.7799
$ws.Range("M138").Value = -1784.4614
$ws.Range("N138").Value = -16440.7799

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 11928.318
$ws.Range("I32").Value = 7368.933
$ws.Range("K32").Value = 7368.933
$ws.Range("M32").Value = -7081.933
# Row 43
$ws.Range("H43").Value = 35000
$ws.Range("J43").Value = 35000
$ws.Range("L43").Value = 35000
$ws.Range("N43").Value = -35626
# Row 45
$ws.Range("H45").Value = 1583.2307
$ws.Range("I45").Value = 1217.8
$ws.Range("J45").Value = 1811.625
$ws.Range("K45").Value = 1217.8
$ws.Range("L45").Value = 1811.625
$ws.Range("M45").Value = -840.8
$ws.Range("N45").Value = -2565.625
# Row 61
$ws.Range("H61").Value = 43696.523
$ws.Range("I61").Value = 51102.5
$ws.Range("K61").Value = 51102.5
$ws.Range("M61").Value = -50890.5
# Row 74
$ws.Range("H74").Value = 965.9375
$ws.Range("I74").Value = 559.4138
$ws.Range("J74").Value = 4895.6665
$ws.Range("K74").Value = 559.4138
$ws.Range("L74").Value = 4895.6665
$ws.Range("M74").Value = 314.5862
$ws.Range("N74").Value = -6643.6665
# Row 77
$ws.Range("H77").Value = 965.9375
$ws.Range("I77").Value = 559.4138
$ws.Range("J77").Value = 4895.6665
$ws.Range("K77").Value = 2797.069
$ws.Range("L77").Value = 24478.3325
$ws.Range("M77").Value = 1570.931
$ws.Range("N77").Value = -33214.3325
# Row 102
$ws.Range("H102").Value = 1277.5
$ws.Range("I102").Value = 1259.3334
$ws.Range("J102").Value = 1550
$ws.Range("K102").Value = 1259.3334
$ws.Range("L102").Value = 1550
$ws.Range("M102").Value = 362.6666
$ws.Range("N102").Value = -4794
# Row 132
$ws.Range("H132").Value = 2346.2258
$ws.Range("I132").Value = 1811.48
$ws.Range("K132").Value = 5434.440000000001
$ws.Range("M132").Value = -2904.440000000001
# Row 136
$ws.Range("H136").Value = 43696.523
$ws.Range("I136").Value = 51102.5
$ws.Range("K136").Value = 153307.5
$ws.Range("M136").Value = -150757.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").ClearContents()
# Row 134
$ws.Range("H134").Value = 5987.8623
$ws.Range("I134").Value = 6875.826
$ws.Range("K134").Value = 20627.478
$ws.Range("M134").Value = -18092.478
# Row 135
$ws.Range("H135").Value = 57780
$ws.Range("J135").Value = 57780
$ws.Range("L135").Value = 57780
$ws.Range("N135").Value = -67920

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 823.1
$ws.Range("I16").Value = 778.875
$ws.Range("K16").Value = 778.875
$ws.Range("M16").Value = -491.875
# Row 58
$ws.Range("H58").Value = 3625524.8
$ws.Range("I58").Value = 5437210.5
$ws.Range("J58").Value = 2153.25
$ws.Range("K58").Value = 5437210.5
$ws.Range("L58").Value = 2153.25
$ws.Range("M58").Value = -5437007.5
$ws.Range("N58").Value = -2559.25
# Row 99
$ws.Range("H99").Value = 590670
$ws.Range("I99").Value = 1002473.8
$ws.Range("K99").Value = 1002473.8
$ws.Range("M99").Value = -1000975.8
# Row 113
$ws.Range("H113").Value = 823.1
$ws.Range("I113").Value = 778.875
$ws.Range("K113").Value = 778.875
$ws.Range("M113").Value = 1391.125
# Row 126
$ws.Range("H126").Value = 590670
$ws.Range("I126").Value = 1002473.8
$ws.Range("K126").Value = 3007421.4
$ws.Range("M126").Value = -3004951.4
# Row 132
$ws.Range("H132").Value = 1918.2142
$ws.Range("I132").Value = 1394.3
$ws.Range("K132").Value = 4182.9
$ws.Range("M132").Value = -1652.9
# Row 136
$ws.Range("H136").Value = 3625524.8
$ws.Range("I136").Value = 5437210.5
$ws.Range("J136").Value = 2153.25
$ws.Range("K136").Value = 16311631.5
$ws.Range("L136").Value = 6459.75
$ws.Range("M136").Value = -16309081.5
$ws.Range("N136").Value = -11559.75

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 82
$ws.Range("H82").Value = 3000
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
# Row 85
$ws.Range("H85").Value = 3000
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
# Row 107
$ws.Range("H107").Value = 1144.7222
$ws.Range("J107").Value = 1413.3077
$ws.Range("L107").Value = 4239.9231
$ws.Range("N107").Value = -8079.9231
# Row 122
$ws.Range("H122").Value = 813.6667
$ws.Range("J122").Value = 881.2857
$ws.Range("L122").Value = 7931.571300000001
$ws.Range("N122").Value = -12831.5713
# Row 139
$ws.Range("H139").Value = 4965.5
$ws.Range("I139").Value = 5075.3335
$ws.Range("J139").Value = 2000
$ws.Range("K139").Value = 15226.0005
$ws.Range("L139").Value = 6000
$ws.Range("M139").Value = -10086.0005
$ws.Range("N139").Value = -16280

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 5999.5
$ws.Range("I102").Value = 5999.5
$ws.Range("K102").Value = 5999.5
$ws.Range("M102").Value = -4377.5
# Row 113
$ws.Range("H113").Value = 1506.2307
$ws.Range("I113").Value = 1142.2858
$ws.Range("K113").Value = 1142.2858
$ws.Range("M113").Value = 1027.7142
# Row 131
$ws.Range("H131").Value = 35555
$ws.Range("J131").Value = 35555
$ws.Range("L131").Value = 35555
$ws.Range("N131").Value = -45635
# Row 132
$ws.Range("H132").Value = 1284204.4
$ws.Range("I132").Value = 1426437.2
$ws.Range("K132").Value = 4279311.6
$ws.Range("M132").Value = -4276781.6
# Row 133
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1128.3334
$ws.Range("I22").Value = 620.7273
$ws.Range("K22").Value = 620.7273
$ws.Range("M22").Value = -325.7273
# Row 27
$ws.Range("H27").Value = 1128.3334
$ws.Range("I27").Value = 620.7273
$ws.Range("K27").Value = 620.7273
$ws.Range("M27").Value = -513.7273
# Row 40
$ws.Range("H40").Value = 9771.556
$ws.Range("I40").Value = 11117.615
$ws.Range("J40").Value = 6271.8
$ws.Range("K40").Value = 11117.615
$ws.Range("L40").Value = 6271.8
$ws.Range("M40").Value = -10981.615
$ws.Range("N40").Value = -6543.8
# Row 46
$ws.Range("H46").Value = 1759.8572
$ws.Range("I46").Value = 1408.6666
$ws.Range("J46").Value = 2023.25
$ws.Range("K46").Value = 1408.6666
$ws.Range("L46").Value = 2023.25
$ws.Range("M46").Value = -1220.6666
$ws.Range("N46").Value = -2399.25
# Row 132
$ws.Range("H132").Value = 2694.182
$ws.Range("I132").Value = 1983.6666
$ws.Range("K132").Value = 5950.9998
$ws.Range("M132").Value = -3420.9998

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 3599.4443
$ws.Range("I81").Value = 1499
$ws.Range("J81").Value = 6225
$ws.Range("K81").Value = 2998
$ws.Range("L81").Value = 12450
$ws.Range("M81").Value = -1937
$ws.Range("N81").Value = -14572
# Row 84
$ws.Range("H84").Value = 3599.4443
$ws.Range("I84").Value = 1499
$ws.Range("J84").Value = 6225
$ws.Range("K84").Value = 14990
$ws.Range("L84").Value = 62250
$ws.Range("M84").Value = -9686
$ws.Range("N84").Value = -72858
# Row 100
$ws.Range("H100").Value = 1034.2858
$ws.Range("I100").Value = 958
$ws.Range("J100").Value = 1225
$ws.Range("K100").Value = 1916
$ws.Range("L100").Value = 2450
$ws.Range("M100").Value = -1375
$ws.Range("N100").Value = -3532
# Row 113
$ws.Range("H113").Value = 592
$ws.Range("I113").Value = 378
$ws.Range("K113").Value = 1134
$ws.Range("M113").Value = 1036
# Row 126
$ws.Range("H126").Value = 5811.52
$ws.Range("I126").Value = 7099.6665
$ws.Range("K126").Value = 21298.9995
$ws.Range("M126").Value = -18828.9995
# Row 131
$ws.Range("H131").Value = 29999.666
$ws.Range("J131").Value = 29999.666
$ws.Range("L131").Value = 29999.666
$ws.Range("N131").Value = -40079.666
# Row 132
$ws.Range("H132").Value = 1204.9512
$ws.Range("I132").Value = 986.8857400000001
$ws.Range("J132").Value = 2477
$ws.Range("K132").Value = 2960.65722
$ws.Range("L132").Value = 7431
$ws.Range("M132").Value = -430.6572200000001
$ws.Range("N132").Value = -12491
# Row 136
$ws.Range("H136").Value = 22223652
$ws.Range("I136").Value = 30864906
$ws.Range("K136").Value = 92594718
$ws.Range("M136").Value = -92592168
